# Applies the "Oppdatert use case diagram" edit:
#  1. Rewrite the "For å gjøre ... har " lead-in sentence.
#  2. Rewrite the tail of that same paragraph (after "vi") with the new
#     use-case-flow sentence, and move the _GoBack bookmark there (right
#     after the new text, at the end of the paragraph).
#  3. Delete the whole following paragraph ("Dette gjør arbeidet på
#     utviklingssiden ... minimalistisk.") entirely.
#  4. The paragraph that used to hold the _GoBack bookmark (just before
#     the sectPr) becomes a plain empty paragraph.

$d = $word.ActiveDocument

# --- 1: lead-in sentence -----------------------------------------------
$d.Content.Find.Execute(
    "For å gjøre det så enkelt som mulig har ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "For å gjøre oppgaven så enkelt gjennomførbar som mulig har ", 2)

# --- 2: tail sentence, replacing the "one page" reasoning --------------
$d.Content.Find.Execute(
    " valgt å plassere alt på én side. ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    " kun noen få sider som brukeren må igjennom. Det første er innlogging, deretter en regitreringsside der brukeren vil få velge hvilket rom de skal bruke.",
    2)

# Locate that paragraph again (its text has just changed) so we can
# anchor the bookmark move precisely at its end (before the paragraph
# mark).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*brukeren vil få velge hvilket rom de skal bruke.*") {
        $target = $para
        break
    }
}

# Remove the old _GoBack bookmark (currently sitting alone in the last
# paragraph of the document, right before the sectPr).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create _GoBack collapsed at the end of $target's text (i.e. right
# before that paragraph's mark). A zero-length Range can't be handed to
# Bookmarks.Add directly, so insert a throwaway marker character, wrap
# the bookmark around it, then delete the marker — leaving the
# bookmark collapsed exactly where the marker was.
$anchor = $target.Range.Duplicate
$anchor.Start = $anchor.End - 1
$anchor.End = $anchor.Start
$anchor.InsertAfter([char]1)
$d.Bookmarks.Add("_GoBack", $anchor)
$markerRange = $d.Bookmarks("_GoBack").Range
$markerRange.Text = ""

# --- 3: delete the following paragraph entirely -------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Dette gjør arbeidet på utviklingssiden*") {
        $para.Range.Delete()
        break
    }
}
